$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = 'ACFT; DITCH; TREE; LOST CONTROL'
$ws.Range('D3').Value = 'TAKEOFF; ENGINE; WING; FUEL TANK; SUMPS; PREFLIGHT'
$ws.Range('D4').Value = 'HELICOPTER; SLING LOAD; LOAD; TREES; PREFLIGHT'
$ws.Range('D5').Value = 'NOSEWHEEL STEERING; BRAKES; AIRCRAFT; FENCE; CIRCUIT BREAKER; HYDRAULIC PUMP'
$ws.Range('D6').Value = 'FORWARD CARGO DOOR; AIRCRAFT; OBJECTS; WARNING LIGHT'
$ws.Range('D7').Value = 'MR. TIMOTHY ALLEN WELLS; PILOT; BELL; HELICOPTER; BHT-47-G5; N4754R'
$ws.Range('D8').Value = 'LANDING; ONE; ENGINE; ICE; AUXILIARY FUEL SYSTEM'
$ws.Range('D9').Value = 'PILOT; ALTITUDE; ICING; LOST CONTROL; ALTIMETER; ICING'
$ws.Range('D10').Value = 'AIRCRAFT; RAY AIRPORT; TAKEOFF; AIRCRAFT; ENGINE; PROBLEMS'
$ws.Range('D11').Value = 'FUEL; FERRY FLIGHT; LEFT PONTOON; LANDING; POND; FUEL CAP; FUEL'
$ws.Range('D12').Value = 'PILOT; TAKEOFF; ENGINE; POWER'
$ws.Range('D13').Value = 'COWLING; TAKEOFF; COWLING; LAND; WINDSHIELD; STABILIZER'
$ws.Range('D14').Value = 'AIRCRAFT; DESCENT; FLIGHT MANUAL; TRUCK; WINDSHIELD; GROUND'
$ws.Range('D15').Value = 'NR2 ENGINE; ENGINE FIRE; PASSENGERS; GATE; FUEL LEVER; SMOKE'
$ws.Range('D16').Value = 'OIL TEMP; GEAR; OIL'
$ws.Range('D17').Value = 'BATTERY COMPARTMENT DOOR; ANTENNA; WINDSHIELD'
$ws.Range('D18').Value = 'NOSE GEAR; PIN; MAINTENANCE'
$ws.Range('D19').Value = 'STAND; BAGGAGE CART; BRAKES; AIRCRAFT; TEAR; FUSELAGE'
$ws.Range('D20').Value = 'NORTHWEST AIRLINES; DC-10; BOEING; GATE B-52; CONTINENTAL AIRLINES'
$ws.Range('D21').Value = 'ENGINE; OIL SMELL; AIR INTAKE PLUGS; FURROWS'
$ws.Range('D22').Value = 'PILOT; ACCESS DOOR; DOOR; FUSELAGE'
$ws.Range('D23').Value = 'RUNWAY; LANDING; AIRPLANE; CARPETING; RUDDER PEDAL'
$ws.Range('D24').Value = 'TAXIWAY; TAXIWAY LIGHT; PROP'
$ws.Range('D25').Value = 'FUEL CAP; RIGHT TANK; GAS; CAP; FLIGHT; WING'
$ws.Range('D26').Value = 'JULY 13, 2005; 1535 MST; RAYTHEON AIRCRAFT CO; B36TC; BONANZA; N3042V; HAVENS LEASING'
$ws.Range('D27').Value = 'ENGINE; POWER; TAKEOFF; AIRCRAFT; LANDING; PILOT; TANKS; WATER; FUEL'
$ws.Range('D28').Value = 'ENGINE; TAKEOFF; ROAD; ICE; CARBURETOR BOWL'
$ws.Range('D29').Value = 'CLIMBOUT; TOW RELEASE; RIGHT AILERON'
$ws.Range('D30').Value = 'AIRCRAFT; MAINTENANCE; GEAR; NOSE BLOCK'
$ws.Range('D31').Value = 'TRAILER; TIE DOWN; SKID; PILOT'
$ws.Range('D32').Value = 'AIRCRAFT; LIFTOFF; RUNWAY; SNOW BANK; FROST; AIRFRAME'
$ws.Range('D33').Value = 'ENGINE; POWER FAILURE; AEROBATICS; FIELD; PILOT; FUEL TANKS'
$ws.Range('D34').Value = 'TAKEOFF; BAGGAGE DOOR'
$ws.Range('D35').Value = 'PILOT; DITCH; TAXIWAY'
$ws.Range('D36').Value = 'TAKEOFF; TREE; MAG SWITCH; ONE MAG'
$ws.Range('D37').Value = 'PILOT; FRONT RANGE AIRPORT; HOURS OF DARKNESS; START AND RUN'
$ws.Range('D38').Value = 'TAKEOFF; PIN; NOSE GEAR'
$ws.Range('D39').Value = 'SIKORSKY; S-70A; N160LA; VIBRATION; APU DOOR; ROTOR BLADE'
$ws.Range('D40').Value = 'PILOT; AIRCRAFT; CARGO TRIP; TAXIWAY'
$ws.Range('D41').Value = 'LANDING; SEATS; OXYGEN MASKS; PASSENGERS'
$ws.Range('D42').Value = 'CARGO DOOR; TAKEOFF; MR.BOWEN; RUNWAY CONDITIONS; STEVEN''S VILLAGE'
$ws.Range('D43').Value = 'TIEDOWN CHAINS; PILOT; RAMP PERSON; SKIDS; WIND'
$ws.Range('D44').Value = 'BAGGAGE CART; AIRCRAFT; JET BLAST; BRAKES; CART'
$ws.Range('D45').Value = 'GLIDER; ALTITUDE'
$ws.Range('D46').Value = 'HELICOPTER; ENGINE; POWER CHECKS; PILOT; CREWMEMBER'
$ws.Range('D47').Value = 'PILOT; LOST CONTROL; AIRCRAFT; TAKEOFF ROLL; PILOTS; SEAT; POSITION; ACCELERATION'
$ws.Range('D48').Value = 'NOSE STRUT; RETRACTION; AIRPLANE; BELLCRANK'
$ws.Range('D49').Value = 'POWER; LANDING; SLOPE; EVIDENCE; FUEL; SCENE'
$ws.Range('D50').Value = 'LEFT ENGINE; H71; VMC; OIL LEAK; ENGINE CASE'
$ws.Range('D51').Value = 'TAKEOFF; SOD STRIP; NORTH WIND; TREES; STRIP'
$ws.Range('D52').Value = 'CRASH; LANDING; ENGINE FAILURE; TAKEOFF; AIRCRAFT; INSPECTION'
$ws.Range('D53').Value = 'PILOT; CESSNA 207; N1549U; TAKEOFF; BETHEL AIRPORT; ENGINE'
$ws.Range('D54').Value = 'FLAPS; TREE; 3 MILES; AIRPORT'
$ws.Range('D55').Value = 'TCA; ATC COMMUNICATION; PILOT; ALTIMETER'
$ws.Range('D56').Value = 'AIRCRAFT; CRUISE FLIGHT; RPM LOSS; PILOT; PILOT; ALTITUDE'
$ws.Range('D57').Value = 'NOSE BAGGAGE DOOR; TAKEOFF; BAGGAGE; LEFT PROP; SEA; ICE; DOOR'
$ws.Range('D58').Value = 'ENGINE; CLIMBOUT; STREET; SPARK PLUGS; TOLERENCE; PRIMER'
$ws.Range('D59').Value = 'TAXI OUT; CFI;LANDING GEAR; LIGHT; LIGHTS'
$ws.Range('D60').Value = 'POWER; IFR DUAL FLIGHT; POWER POLE; LANDING; CIRCUIT BREAKER; PILOT ERROR'
$ws.Range('D61').Value = 'BANG; SHUDDER; SPARKS; ROTOR; COWLING; COWL'
$ws.Range('D62').Value = 'PROBLEM; CLIMBOUT; LOCKING PINS; SEAT TRACK; SEAT PIN; BELT'
$ws.Range('D63').Value = 'LANDING; POWER LOSS; WATER; FUEL'
$ws.Range('D64').Value = 'AIRPLANE; TAKEOFF ROLL; AILERON-ELEVATOR LOCK PIN'
$ws.Range('D65').Value = 'PILOT; RAIN SHOWER; CLIMBOUT; ENGINE; FIELD; WATER; FUEL'
$ws.Range('D66').Value = 'PILOT; TAKEOFF; SEAT'
$ws.Range('D67').Value = 'MARCH 14, 1995; N7016M; AIRCRAFT; EMERGENCY OFF AIRPORT LANDING'
$ws.Range('D68').Value = 'FUEL EXHAUSTION; SIPHONING; VENT CAP; EMERGENCY LANDING AREA'
$ws.Range('D69').Value = 'AIRCRAFT; AIRCRAFT; BRAKES'
$ws.Range('D70').Value = 'ENGINE START; PILOTS; SEAT; THROTTLE; AIRCRAFT; PLANE'
$ws.Range('D71').Value = 'LIFTOFF; PILOT; RIGHT FUEL CAP; RUNWAY; RUNWAY'
$ws.Range('D72').Value = 'AIRCRAFT; BARTLESVILLE MUNICIPAL AIRPORT; BVO; BARTLESVILLE, OK; 1147 AM LOCAL TIME'
$ws.Range('D73').Value = 'PILOT; RUNWAY; BRAKE'
$ws.Range('D74').Value = 'WHEELS; LANDINGS; SNOW'
$ws.Range('D75').Value = 'POWER; WATER; FLOATS; LONGERON; WATER; FUEL; SUMP'
$ws.Range('D76').Value = 'TOUCHDOWN; AIRCRAFT; RUNWAY; AIRCRAFT'
$ws.Range('D77').Value = 'COWLING; ENGINE; FLIGHT'
$ws.Range('D78').Value = 'FLIGHT; NOTH PLATTER, NE; IOWA CITY, IA; 487 STATUTE MILES; AIRCRAFT; POWER; GRINNELL'
$ws.Range('D79').Value = 'TIE DOWN STRAP; LEFT SKID; TRAILER'
$ws.Range('D80').Value = 'PILOT IN COMMAND; PIC; FUEL TANK; CAPS; DEPARTURE'
$ws.Range('D81').Value = 'FUEL CONTAMINATION; AIRFRAME; FUEL FILTER; FUEL FILTER; ENGINE; FUEL PUMP FILTER'
$ws.Range('D82').Value = 'DEMO FLIGHT; ROD; HYRDO LOCK; NR5 CYL; STARTUP; COUNTER WT; SHAFT'
$ws.Range('D83').Value = 'AIRCRAFT; N106DA; PILOT; BRAKE PROBLEMS; AIRCRAFT'
$ws.Range('D84').Value = 'N759TY; CE-182; HOT SPRINGS, ARKANSAS; FAYETTEVILLE, TENNESSEE; AIRCRAFT; ENGINE'
$ws.Range('D85').Value = 'CLIMBOUT; NOSE; BAGGAGE DOOR; POWER; OCEAN; TEST FLIGHT; STORAGE'
$ws.Range('D86').Value = 'ELEVATOR CONTROL; HARNESS; ELEVATOR CABLE; STICK BOOT; FLOOR; FIELD'
$ws.Range('D87').Value = 'POWER; STRIP;LEFT FUEL CAP; GRAVITY FUEL FLOW'
$ws.Range('D88').Value = 'ACFT DISPATCHER; PILOT; PILOT'
$ws.Range('D89').Value = 'HELICOPTER; HELIPAD; RPM CONTROL; FLY POSITION; TAKEOFF'
$ws.Range('D90').Value = 'TUBE; AERIAL BANNERS; ROOF; HOUSE; OPERATOR'
$ws.Range('D91').Value = 'MR.KADERA; FIELD; HIGHWAY 93; THREE; EAST; SUNMER, IOWA'
$ws.Range('D92').Value = 'SEPTEMBER 23, 1999; 1900 HOURS PACIFIC DAYLIGHT TIME; CIERNIA GLASAIR III; N153JC'
$ws.Range('D93').Value = 'NOSE; BAGGAGE; CARGO; WHEEL CHOCK; FLIGHT; NOSE LANDING GEAR'
$ws.Range('D94').Value = 'EXPLOSION; LOGS; COMPRESSOR STALL; NR1 ENGINE; LOGS; SE PROCEDURES'
$ws.Range('D95').Value = 'FEBRUARY 19, 2000; 1825 EASTERN STANDARD TIME; BEECH 1900D; N81SK; SKYWAY AIRLINE; FLIGHT'
$ws.Range('D96').Value = 'GROUND LOOPED; LANDING; TAILWHEEL; RIGHT'
$ws.Range('D97').Value = 'ENGINE; CLIMBOUT; CRASH LANDED; TREES; FUEL TANK; TRIP; FUEL'
$ws.Range('D98').Value = 'WITNESSES; DEPARTURE; RUNWAY; AIRPLANE'
$ws.Range('D99').Value = 'POP; TAKEOFF ROLL; LOST RUDDER CONTROL; RWY; CORNFIELD; TAILWHEEL'
$ws.Range('D100').Value = 'CANOPY; CLIMBOUT; GROUND COLLISION; DEFECTS; PILOT ERROR'
$ws.Range('D101').Value = 'PILOT; OIL FILLER CAP; OIL; ENGINE'
